# Weekly update: three new daily price-report rows are added to the
# "Ciruela" (plum) price sheet for Vega Monumental Concepción.
#
#  - A brand-new row is inserted at row 83 (newest date, 2023-03-22).
#  - Two brand-new rows are inserted at rows 120-121 (2023-03-23), which
#    pushes all the older historical rows down accordingly.
#
# All pre-existing rows keep their original values; they are simply
# shifted down by the corresponding number of inserted rows.

function Set-RowValues($Sheet, $Row, $Values) {
    $n = $Values.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCell = $Sheet.Cells.Item($Row, 1)
    $endCell = $Sheet.Cells.Item($Row, $n)
    $range = $Sheet.Range($startCell, $endCell)
    $range.Value2 = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row at position 83 -------------------------------
$ws.Rows.Item(83).Insert()

Set-RowValues $ws 83 @(
    11,
    'Vega Monumental Concepción',
    'Bíobío',
    45007,
    8,
    'Fruta',
    100103,
    'Frutos de hueso (carozo)',
    100103002,
    'Ciruela',
    'Angeleno',
    'Primera',
    270,
    10000,
    11000,
    10556,
    '$/bandeja 18 kilos granel',
    "Región de O'Higgins",
    586,
    18
)

# --- Insert the two new rows at positions 120-121 ---------------------
$ws.Range("A120:A121").EntireRow.Insert()

Set-RowValues $ws 120 @(
    11,
    'Vega Monumental Concepción',
    'Bíobío',
    45008,
    8,
    'Fruta',
    100103,
    'Frutos de hueso (carozo)',
    100103002,
    'Ciruela',
    'Angeleno',
    'Especial',
    220,
    10000,
    11000,
    10455,
    '$/bandeja 18 kilos granel',
    'Provincia de Curicó',
    581,
    18
)

Set-RowValues $ws 121 @(
    11,
    'Vega Monumental Concepción',
    'Bíobío',
    45008,
    8,
    'Fruta',
    100103,
    'Frutos de hueso (carozo)',
    100103002,
    'Ciruela',
    'Angeleno',
    'Primera',
    220,
    8500,
    9000,
    8727,
    '$/bandeja 18 kilos granel',
    'Provincia de Curicó',
    485,
    18
)
